$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MayRaw")

# Header row
$ws.Cells.Item(1,1).Value = "Library"
$ws.Cells.Item(1,2).Value = "Items owned by this library checked out at this library this month"
$ws.Cells.Item(1,3).Value = "Items owned by other libraries checked out at this library this month"
$ws.Cells.Item(1,4).Value = "Total circulation this month"

# Data rows 2-54
$ws.Cells.Item(2,1).Value = "Atchison Public Library"
$ws.Cells.Item(2,2).Value = 4067
$ws.Cells.Item(2,3).Value = 1390
$ws.Cells.Item(2,4).Value = 5457

$ws.Cells.Item(3,1).Value = "Baldwin City Public Library"
$ws.Cells.Item(3,2).Value = 2599
$ws.Cells.Item(3,3).Value = 524
$ws.Cells.Item(3,4).Value = 3123

$ws.Cells.Item(4,1).Value = "Basehor Community Library"
$ws.Cells.Item(4,2).Value = 7350
$ws.Cells.Item(4,3).Value = 1129
$ws.Cells.Item(4,4).Value = 8479

$ws.Cells.Item(5,1).Value = "Bern Community Library"
$ws.Cells.Item(5,2).Value = 88
$ws.Cells.Item(5,3).Value = 30
$ws.Cells.Item(5,4).Value = 118

$ws.Cells.Item(6,1).Value = "Bonner Springs City Library"
$ws.Cells.Item(6,2).Value = 4809
$ws.Cells.Item(6,3).Value = 1132
$ws.Cells.Item(6,4).Value = 5941

$ws.Cells.Item(7,1).Value = "Burlingame Community Library"
$ws.Cells.Item(7,2).Value = 437
$ws.Cells.Item(7,3).Value = 211
$ws.Cells.Item(7,4).Value = 648

$ws.Cells.Item(8,1).Value = "Carbondale City Library"
$ws.Cells.Item(8,2).Value = 419
$ws.Cells.Item(8,3).Value = 117
$ws.Cells.Item(8,4).Value = 536

$ws.Cells.Item(9,1).Value = "Centralia Community Library"
$ws.Cells.Item(9,2).Value = 220
$ws.Cells.Item(9,3).Value = 46
$ws.Cells.Item(9,4).Value = 266

$ws.Cells.Item(10,1).Value = "Corning City Library"
$ws.Cells.Item(10,2).Value = 23
$ws.Cells.Item(10,4).Value = 23

$ws.Cells.Item(11,1).Value = "Digital Content"

$ws.Cells.Item(12,1).Value = "Doniphan County Library - Elwood"
$ws.Cells.Item(12,2).Value = 75
$ws.Cells.Item(12,3).Value = 20
$ws.Cells.Item(12,4).Value = 95

$ws.Cells.Item(13,1).Value = "Doniphan County Library - Highland"
$ws.Cells.Item(13,2).Value = 144
$ws.Cells.Item(13,3).Value = 77
$ws.Cells.Item(13,4).Value = 221

$ws.Cells.Item(14,1).Value = "Doniphan County Library - Troy"
$ws.Cells.Item(14,2).Value = 461
$ws.Cells.Item(14,3).Value = 254
$ws.Cells.Item(14,4).Value = 715

$ws.Cells.Item(15,1).Value = "Doniphan County Library - Wathena"
$ws.Cells.Item(15,2).Value = 301
$ws.Cells.Item(15,3).Value = 56
$ws.Cells.Item(15,4).Value = 357

$ws.Cells.Item(16,1).Value = "Effingham Community Library"
$ws.Cells.Item(16,2).Value = 253
$ws.Cells.Item(16,3).Value = 66
$ws.Cells.Item(16,4).Value = 319

$ws.Cells.Item(17,1).Value = "Eudora Community Library"
$ws.Cells.Item(17,2).Value = 1475
$ws.Cells.Item(17,3).Value = 620
$ws.Cells.Item(17,4).Value = 2095

$ws.Cells.Item(18,1).Value = "Everest, Barnes Reading Room"
$ws.Cells.Item(18,2).Value = 78
$ws.Cells.Item(18,3).Value = 71
$ws.Cells.Item(18,4).Value = 149

$ws.Cells.Item(19,1).Value = "Hiawatha, Morrill Public Library"
$ws.Cells.Item(19,2).Value = 1575
$ws.Cells.Item(19,3).Value = 557
$ws.Cells.Item(19,4).Value = 2132

$ws.Cells.Item(20,1).Value = "Highland Community College"
$ws.Cells.Item(20,2).Value = 48
$ws.Cells.Item(20,3).Value = 33
$ws.Cells.Item(20,4).Value = 81

$ws.Cells.Item(21,1).Value = "Holton, Beck-Bookman Library"
$ws.Cells.Item(21,2).Value = 1680
$ws.Cells.Item(21,3).Value = 507
$ws.Cells.Item(21,4).Value = 2187

$ws.Cells.Item(22,1).Value = "Horton Public Library"
$ws.Cells.Item(22,2).Value = 198
$ws.Cells.Item(22,3).Value = 90
$ws.Cells.Item(22,4).Value = 288

$ws.Cells.Item(23,1).Value = "Lansing Community Library"
$ws.Cells.Item(23,2).Value = 2004
$ws.Cells.Item(23,3).Value = 616
$ws.Cells.Item(23,4).Value = 2620

$ws.Cells.Item(24,1).Value = "Leavenworth Public Library"
$ws.Cells.Item(24,2).Value = 8467
$ws.Cells.Item(24,3).Value = 1710
$ws.Cells.Item(24,4).Value = 10177

$ws.Cells.Item(25,1).Value = "Linwood Community Library"
$ws.Cells.Item(25,2).Value = 563
$ws.Cells.Item(25,3).Value = 203
$ws.Cells.Item(25,4).Value = 766

$ws.Cells.Item(26,1).Value = "Louisburg Library"

$ws.Cells.Item(27,1).Value = "Lyndon Carnegie Library"
$ws.Cells.Item(27,2).Value = 253
$ws.Cells.Item(27,3).Value = 201
$ws.Cells.Item(27,4).Value = 454

$ws.Cells.Item(28,1).Value = "McLouth Public Library"
$ws.Cells.Item(28,2).Value = 255
$ws.Cells.Item(28,3).Value = 94
$ws.Cells.Item(28,4).Value = 349

$ws.Cells.Item(29,1).Value = "Meriden-Ozawkie Public Library"
$ws.Cells.Item(29,2).Value = 1430
$ws.Cells.Item(29,3).Value = 571
$ws.Cells.Item(29,4).Value = 2001

$ws.Cells.Item(30,1).Value = "Northeast Kansas Library System"
$ws.Cells.Item(30,2).Value = 14
$ws.Cells.Item(30,3).Value = 38
$ws.Cells.Item(30,4).Value = 52

$ws.Cells.Item(31,1).Value = "Nortonville Public Library"
$ws.Cells.Item(31,2).Value = 268
$ws.Cells.Item(31,3).Value = 78
$ws.Cells.Item(31,4).Value = 346

$ws.Cells.Item(32,1).Value = "Osage City Library"
$ws.Cells.Item(32,2).Value = 1277
$ws.Cells.Item(32,3).Value = 371
$ws.Cells.Item(32,4).Value = 1648

$ws.Cells.Item(33,1).Value = "Osawatomie Public Library"
$ws.Cells.Item(33,2).Value = 858
$ws.Cells.Item(33,3).Value = 443
$ws.Cells.Item(33,4).Value = 1301

$ws.Cells.Item(34,1).Value = "Oskaloosa Public Library"
$ws.Cells.Item(34,2).Value = 514
$ws.Cells.Item(34,3).Value = 193
$ws.Cells.Item(34,4).Value = 707

$ws.Cells.Item(35,1).Value = "Ottawa Library"
$ws.Cells.Item(35,2).Value = 5769
$ws.Cells.Item(35,3).Value = 787
$ws.Cells.Item(35,4).Value = 6556

$ws.Cells.Item(36,1).Value = "Overbrook Public Library"
$ws.Cells.Item(36,2).Value = 685
$ws.Cells.Item(36,3).Value = 172
$ws.Cells.Item(36,4).Value = 857

$ws.Cells.Item(37,1).Value = "Paola Free Library"
$ws.Cells.Item(37,2).Value = 3134
$ws.Cells.Item(37,3).Value = 485
$ws.Cells.Item(37,4).Value = 3619

$ws.Cells.Item(38,1).Value = "Perry-Lecompton Community Library"
$ws.Cells.Item(38,2).Value = 172
$ws.Cells.Item(38,3).Value = 16
$ws.Cells.Item(38,4).Value = 188

$ws.Cells.Item(39,1).Value = "Pomona Community Library"
$ws.Cells.Item(39,2).Value = 67
$ws.Cells.Item(39,3).Value = 60
$ws.Cells.Item(39,4).Value = 127

$ws.Cells.Item(40,1).Value = "Prairie Hills Schools - Axtell Public School"
$ws.Cells.Item(40,2).Value = 384
$ws.Cells.Item(40,3).Value = 15
$ws.Cells.Item(40,4).Value = 399

$ws.Cells.Item(41,1).Value = "Prairie Hills Schools - Sabetha Elementary School"
$ws.Cells.Item(41,2).Value = 2224
$ws.Cells.Item(41,3).Value = 72
$ws.Cells.Item(41,4).Value = 2296

$ws.Cells.Item(42,1).Value = "Prairie Hills Schools - Sabetha High School"
$ws.Cells.Item(42,2).Value = 31
$ws.Cells.Item(42,3).Value = 10
$ws.Cells.Item(42,4).Value = 41

$ws.Cells.Item(43,1).Value = "Prairie Hills Schools - Sabetha Middle School"
$ws.Cells.Item(43,2).Value = 142
$ws.Cells.Item(43,3).Value = 8
$ws.Cells.Item(43,4).Value = 150

$ws.Cells.Item(44,1).Value = "Richmond Public Library"
$ws.Cells.Item(44,2).Value = 343
$ws.Cells.Item(44,3).Value = 75
$ws.Cells.Item(44,4).Value = 418

$ws.Cells.Item(45,1).Value = "Rossville Community Library"
$ws.Cells.Item(45,2).Value = 1364
$ws.Cells.Item(45,3).Value = 606
$ws.Cells.Item(45,4).Value = 1970

$ws.Cells.Item(46,1).Value = "Sabetha, Mary Cotton Library"
$ws.Cells.Item(46,2).Value = 2435
$ws.Cells.Item(46,3).Value = 1122
$ws.Cells.Item(46,4).Value = 3557

$ws.Cells.Item(47,1).Value = "Seneca Free Library"
$ws.Cells.Item(47,2).Value = 1478
$ws.Cells.Item(47,3).Value = 233
$ws.Cells.Item(47,4).Value = 1711

$ws.Cells.Item(48,1).Value = "Silver Lake Library"
$ws.Cells.Item(48,2).Value = 889
$ws.Cells.Item(48,3).Value = 658
$ws.Cells.Item(48,4).Value = 1547

$ws.Cells.Item(49,1).Value = "Tonganoxie Public Library"
$ws.Cells.Item(49,2).Value = 2627
$ws.Cells.Item(49,3).Value = 614
$ws.Cells.Item(49,4).Value = 3241

$ws.Cells.Item(50,1).Value = "Valley Falls, Delaware Township Library"
$ws.Cells.Item(50,2).Value = 483
$ws.Cells.Item(50,3).Value = 163
$ws.Cells.Item(50,4).Value = 646

$ws.Cells.Item(51,1).Value = "Wellsville City Library"
$ws.Cells.Item(51,2).Value = 881
$ws.Cells.Item(51,3).Value = 371
$ws.Cells.Item(51,4).Value = 1252

$ws.Cells.Item(52,1).Value = "Wetmore Public Library"
$ws.Cells.Item(52,2).Value = 87
$ws.Cells.Item(52,3).Value = 163
$ws.Cells.Item(52,4).Value = 250

$ws.Cells.Item(53,1).Value = "Williamsburg Community Library"
$ws.Cells.Item(53,2).Value = 174
$ws.Cells.Item(53,3).Value = 51
$ws.Cells.Item(53,4).Value = 225

$ws.Cells.Item(54,1).Value = "Winchester Public Library"
$ws.Cells.Item(54,2).Value = 361
$ws.Cells.Item(54,3).Value = 365
$ws.Cells.Item(54,4).Value = 726

# Make "May" the active/selected sheet, with B2 as the active cell
# (previously "January" was the tabSelected sheet)
$mayWs = $wb.Worksheets.Item("May")
$mayWs.Activate()
$mayWs.Range("B2").Select()
